# Insert a new employee record as row 2 (right after the header row),
# shifting all existing data rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every existing data row down one slot.
$ws.Rows.Item(2).Insert()

# Fill in the new row's four columns: employeeNumber, fullName, branch, category.
# Branch (C2) is written before fullName (B2) so new shared-string entries are
# appended in the same order as the target workbook (branch text first).
$ws.Range("A2").Value = 1111
$ws.Range("C2").Value = "1111 - sucu"
$ws.Range("B2").Value = "agus millan tst"
$ws.Range("D2").Value = "C-JU - EC Junior"

# Match the author's final cursor position.
$ws.Range("B12").Select()
